# edit.ps1 -- apply the "Cosmos to Culture" -> "Exploring the Universe of
# Chemistry" rewrite described by the supplied diff.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                                      $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "MISSING: $old"
    }
}

# --- Title ---------------------------------------------------------------
Replace-Text "Cosmos to Culture: Cosmic Cycles in Human Phenomena" `
             "Exploring the Universe of Chemistry"

# --- Author (collapses the 3 "Dr" / "." / " Alan Fields" runs into one) --
Replace-Text "Dr. Alan Fields" "Thomas Williamson"

# --- Email address ---------------------------------------------------------
# Original runs: "af8866@emailworld" + "." + "net"
# Target runs:   "thomas" + "." + "williamson@knightsbridge" + "." + "net"
# Only the first run's text actually changes; the trailing "." and "net"
# runs are untouched.
Replace-Text "af8866@emailworld" "thomas.williamson@knightsbridge"

# --- Body paragraph (intro) ------------------------------------------------
Replace-Text "In the kaleidoscopic tapestry of existence, the universe and human culture exhibit a captivating interplay of cycles--an intricate dance between cosmic patterns and terrestrial expressions" `
             "Chemistry, the study of matter, its properties, and its behavior, is an enthralling field that holds the key to understanding the world around us"

Replace-Text " From the eternal rhythm of celestial bodies to the ebb and flow of human emotions, civilizations, and art, the ceaseless recurrence of cycles weaves a rich tapestry of interconnectedness" `
             " This field is essential in uncovering the intricate workings of the universe, from the tiniest atoms to the vast cosmos"

Replace-Text " This essay delves into this cosmic synchronicity, exploring the profound influence of celestial cycles on human phenomena, including art, literature, music, and societal structures" `
             " From the earliest alchemists to modern scientists, the exploration of chemistry has fueled innovation, revolutionized industries, and shaped human history"

Replace-Text "The cosmos, an enigmatic expanse of mysteries, unveils the majesty of repeating celestial patterns" `
             "Chemistry is not merely a collection of theories and formulas; it is an active pursuit of knowledge, a relentless quest to uncover the hidden secrets of nature"

Replace-Text " Day and night, month by month, season to season, the movement of planets, the alignment of constellations, and the dance of galaxies paint a cosmic symphony, influencing the rhythms of life on Earth" `
             " Each experiment conducted, each discovery made, brings us closer to unraveling the profound mysteries of the universe"

Replace-Text " Ancient civilizations looked to the skies for guidance, their rituals and mythologies imbued with celestial symbolism" `
             " Chemistry enables us to understand the composition of materials, their interactions, and their transformations"

Replace-Text " Even today, astrology and other forms of celestial divination bear witness to the enduring allure of the cosmos" `
             " It empowers us to create new substances, design advanced materials, and harness energy in novel ways"

Replace-Text "Furthermore, the cycles of human life mirror the rhythms of the universe" `
             "The study of chemistry is a gateway to a world of endless possibilities"

Replace-Text " Birth, growth, decay, and renewal--the stages of life--resemble the cosmic cycle of expansion and contraction" `
             " It is a journey that takes us from the microscopic realm of atoms and molecules to the macroscopic scale of chemical processes that shape our world"

Replace-Text " The natural world, with its seasonal cycles of growth and dormancy, mirrors the human experience of joy and sorrow, triumph and adversity, love and loss" `
             " Chemistry offers a lens through which we can comprehend the natural world, the human body, and the universe itself"

Replace-Text " These parallels invite introspection, urging us to seek harmony with the cycles of nature and the universe" `
             " It is a field that fosters curiosity, ignites the imagination, and empowers us to solve complex problems"

# --- Summary heading paragraph stays "Summary" (unchanged) ----------------

# --- Summary body -----------------------------------------------------------
Replace-Text "In essence, the cosmic cycles that govern the universe resonate deeply with human existence, inspiring art, shaping cultural narratives, and mirroring the rhythms of life" `
             "Chemistry is the science that studies the properties, behavior, and interactions of matter"

Replace-Text " The recurrence of celestial patterns, such as the lunar cycle and the solstices, finds expression in human rituals, traditions, and mythology" `
             " It plays a fundamental role in our understanding of the world, from the smallest particles to the vastness of the universe"

Replace-Text " The cyclical nature of life--birth, growth, decay, and renewal--mirrors the cosmic dance of expansion and contraction" `
             " Chemistry has led to groundbreaking discoveries, shaping industries, and revolutionizing our understanding of nature"

Replace-Text " This essay highlights the profound interconnectedness between the cosmos and human culture, inviting contemplation of our place within the grand tapestry of the universe" `
             " Through its exploration, we continue to uncover the secrets of the universe and enhance our ability to harness its potential to improve lives and solve global challenges"

# --- Add the new trailing empty paragraph at the end of the document ------
$d.Paragraphs.Add() | Out-Null

Write-Output "done"
